{"js": "const body = context.document.body;\nconst replacements = [\n  [\"2023-09-01 Friday\", \"2023-09-02 Saturday\"],\n  [\"25\u00d776=\", \"36\u00d799=\"],\n  [\"52\u00d715=\", \"61\u00d797=\"],\n  [\"84\u00d725=\", \"82\u00d731=\"],\n  [\"44\u00d728=\", \"13\u00d794=\"],\n  [\"39\u00d740=\", \"50\u00d768=\"],\n  [\"62\u00d773=\", \"71\u00d785=\"],\n  [\"52\u00d716=\", \"57\u00d798=\"],\n  [\"62\u00d746=\", \"87\u00d751=\"],\n  [\"85\u00d738=\", \"75\u00d784=\"],\n  [\"50\u00d742=\", \"74\u00d726=\"],\n  [\"28\u00d756=\", \"44\u00d779=\"],\n  [\"66\u00d735=\", \"23\u00d764=\"],\n  [\"37\u00d762=\", \"12\u00d781=\"],\n  [\"17\u00d740=\", \"44\u00d735=\"],\n  [\"57\u00d786=\", \"50\u00d775=\"],\n  [\"20\u00d720=\", \"21\u00d720=\"],\n  [\"87\u00d777=\", \"59\u00d772=\"],\n  [\"36\u00d720=\", \"72\u00d796=\"],\n  [\"25\u00d749=\", \"56\u00d797=\"],\n  [\"69\u00d759=\", \"63\u00d740=\"],\n  [\"25\u00d775=\", \"36\u00d756=\"],\n  [\"57\u00d785=\", \"94\u00d748=\"],\n  [\"61\u00d778=\", \"21\u00d792=\"],\n  [\"77\u00d712=\", \"75\u00d735=\"],\n  [\"38\u00d776=\", \"82\u00d788=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n$find = $d.Content.Find\n\n$find.Execute('2023-09-01 Friday', $false, $false, $false, $false, $false, $true, 1, $false, '2023-09-02 Saturday', 2) | Out-Null\n$find.Execute('25\u00d776=', $false, $false, $false, $false, $false, $true, 1, $false, '36\u00d799=', 2) | Out-Null\n$find.Execute('52\u00d715=', $false, $false, $false, $false, $false, $true, 1, $false, '61\u00d797=', 2) | Out-Null\n$find.Execute('84\u00d725=', $false, $false, $false, $false, $false, $true, 1, $false, '82\u00d731=', 2) | Out-Null\n$find.Execute('44\u00d728=', $false, $false, $false, $false, $false, $true, 1, $false, '13\u00d794=', 2) | Out-Null\n$find.Execute('39\u00d740=', $false, $false, $false, $false, $false, $true, 1, $false, '50\u00d768=', 2) | Out-Null\n$find.Execute('62\u00d773=', $false, $false, $false, $false, $false, $true, 1, $false, '71\u00d785=', 2) | Out-Null\n$find.Execute('52\u00d716=', $false, $false, $false, $false, $false, $true, 1, $false, '57\u00d798=', 2) | Out-Null\n$find.Execute('62\u00d746=', $false, $false, $false, $false, $false, $true, 1, $false, '87\u00d751=', 2) | Out-Null\n$find.Execute('85\u00d738=', $false, $false, $false, $false, $false, $true, 1, $false, '75\u00d784=', 2) | Out-Null\n$find.Execute('50\u00d742=', $false, $false, $false, $false, $false, $true, 1, $false, '74\u00d726=', 2) | Out-Null\n$find.Execute('28\u00d756=', $false, $false, $false, $false, $false, $true, 1, $false, '44\u00d779=', 2) | Out-Null\n$find.Execute('66\u00d735=', $false, $false, $false, $false, $false, $true, 1, $false, '23\u00d764=', 2) | Out-Null\n$find.Execute('37\u00d762=', $false, $false, $false, $false, $false, $true, 1, $false, '12\u00d781=', 2) | Out-Null\n$find.Execute('17\u00d740=', $false, $false, $false, $false, $false, $true, 1, $false, '44\u00d735=', 2) | Out-Null\n$find.Execute('57\u00d786=', $false, $false, $false, $false, $false, $true, 1, $false, '50\u00d775=', 2) | Out-Null\n$find.Execute('20\u00d720=', $false, $false, $false, $false, $false, $true, 1, $false, '21\u00d720=', 2) | Out-Null\n$find.Execute('87\u00d777=', $false, $false, $false, $false, $false, $true, 1, $false, '59\u00d772=', 2) | Out-Null\n$find.Execute('36\u00d720=', $false, $false, $false, $false, $false, $true, 1, $false, '72\u00d796=', 2) | Out-Null\n$find.Execute('25\u00d749=', $false, $false, $false, $false, $false, $true, 1, $false, '56\u00d797=', 2) | Out-Null\n$find.Execute('69\u00d759=', $false, $false, $false, $false, $false, $true, 1, $false, '63\u00d740=', 2) | Out-Null\n$find.Execute('25\u00d775=', $false, $false, $false, $false, $false, $true, 1, $false, '36\u00d756=', 2) | Out-Null\n$find.Execute('57\u00d785=', $false, $false, $false, $false, $false, $true, 1, $false, '94\u00d748=', 2) | Out-Null\n$find.Execute('61\u00d778=', $false, $false, $false, $false, $false, $true, 1, $false, '21\u00d792=', 2) | Out-Null\n$find.Execute('77\u00d712=', $false, $false, $false, $false, $false, $true, 1, $false, '75\u00d735=', 2) | Out-Null\n$find.Execute('38\u00d776=', $false, $false, $false, $false, $false, $true, 1, $false, '82\u00d788=', 2) | Out-Null\n"}
